# Update CDA Logical model for ST.r2b
# - Bump Version and Date metadata values
# - Insert a new "Jurisdiction" property row (empty value) right after "Contact"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) Update the Version value (row 3, column B)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# 2) Update the Date value (row 8, column B)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# 3) Insert a new row right before the "Description" row (row 11) so the rest of
#    the property table shifts down by one, then fix up the new row's formatting
#    to match its neighbours (Insert() alone creates a distinct style id).
$ws.Rows.Item(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
